$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 500
$ws.Range("K12").Value = 500
$ws.Range("M12").Value = -330

$ws.Range("H40").Value = 1958.4814
$ws.Range("I40").Value = 1884.15
$ws.Range("K40").Value = 1884.15
$ws.Range("M40").Value = -1709.15

$ws.Range("H74").Value = 173609.22
$ws.Range("I74").Value = 173609.22
$ws.Range("K74").Value = 173609.22
$ws.Range("M74").Value = -172673.22

$ws.Range("H77").Value = 173609.22
$ws.Range("I77").Value = 173609.22
$ws.Range("K77").Value = 868046.1
$ws.Range("M77").Value = -863366.1

$ws.Range("H92").Value = 300.7
$ws.Range("I92").Value = 265.2857
$ws.Range("K92").Value = 265.2857
$ws.Range("M92").Value = 982.7143

$ws.Range("H94").Value = 905
$ws.Range("I94").Value = 905
$ws.Range("K94").Value = 905
$ws.Range("M94").Value = -454

$ws.Range("H106").Value = 3596.8
$ws.Range("I106").Value = 3596.8
$ws.Range("K106").Value = 3596.8
$ws.Range("M106").Value = -2965.8

$ws.Range("H107").Value = 599.3333
$ws.Range("I107").Value = 599.5
$ws.Range("J107").Value = 599
$ws.Range("K107").Value = 599.5
$ws.Range("L107").Value = 599
$ws.Range("M107").Value = 1320.5
$ws.Range("N107").Value = -4439

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5084.533
$ws.Range("I61").Value = 1162.6666
$ws.Range("K61").Value = 1162.6666
$ws.Range("M61").Value = -950.6666

$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802

$ws.Range("H136").Value = 5084.533
$ws.Range("I136").Value = 1162.6666
$ws.Range("K136").Value = 3487.9998
$ws.Range("M136").Value = -937.9998000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 9249.75
$ws.Range("I82").Value = 9249.75
$ws.Range("K82").Value = 9249.75
$ws.Range("M82").Value = -8866.75

$ws.Range("H85").Value = 9249.75
$ws.Range("I85").Value = 9249.75
$ws.Range("K85").Value = 9249.75
$ws.Range("M85").Value = -7923.75

$ws.Range("H99").Value = 1782.0834
$ws.Range("I99").Value = 1398.7273
$ws.Range("J99").Value = 5999
$ws.Range("K99").Value = 1398.7273
$ws.Range("L99").Value = 5999
$ws.Range("M99").Value = 99.27269999999999
$ws.Range("N99").Value = -8995

$ws.Range("H107").Value = 2998.3333
$ws.Range("I107").Value = 2998.3333
$ws.Range("K107").Value = 2998.3333
$ws.Range("M107").Value = -1078.3333

$ws.Range("H122").Value = 70387.5
$ws.Range("J122").Value = 70387.5
$ws.Range("L122").Value = 70387.5
$ws.Range("N122").Value = -80187.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1841.6666
$ws.Range("I58").Value = 1758.75
$ws.Range("K58").Value = 1758.75
$ws.Range("M58").Value = -1555.75

$ws.Range("H92").Value = 32200
$ws.Range("J92").Value = 32200
$ws.Range("L92").Value = 32200
$ws.Range("N92").Value = -37192

$ws.Range("H134").Value = 2360.077
$ws.Range("I134").Value = 2831.5557
$ws.Range("J134").Value = 1299.25
$ws.Range("K134").Value = 8494.667099999999
$ws.Range("L134").Value = 3897.75
$ws.Range("M134").Value = -5959.667099999999
$ws.Range("N134").Value = -8967.75

$ws.Range("H136").Value = 1841.6666
$ws.Range("I136").Value = 1758.75
$ws.Range("K136").Value = 5276.25
$ws.Range("M136").Value = -2726.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12222913
$ws.Range("I4").Value = 12222913
$ws.Range("K4").Value = 36668739
$ws.Range("M4").Value = -36668627

$ws.Range("H17").Value = 500
$ws.Range("J17").Value = 500
$ws.Range("L17").Value = 1500
$ws.Range("N17").Value = -1838

$ws.Range("H49").Value = 2321.2856
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 2321.2856
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 6963.8568
$ws.Range("M49").Value = ""
$ws.Range("N49").Value = -7275.8568

$ws.Range("H86").Value = 11500
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = ""

$ws.Range("H89").Value = 11500
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = ""

$ws.Range("H129").Value = 1002740.2
$ws.Range("I129").Value = 1400
$ws.Range("J129").Value = 1253075.2
$ws.Range("K129").Value = 4200
$ws.Range("L129").Value = 3759225.6
$ws.Range("M129").Value = 800
$ws.Range("N129").Value = -3769225.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8413.263000000001
$ws.Range("J70").Value = 9977.223
$ws.Range("L70").Value = 9977.223
$ws.Range("N70").Value = -10517.223

$ws.Range("H73").Value = 8413.263000000001
$ws.Range("J73").Value = 9977.223
$ws.Range("L73").Value = 9977.223
$ws.Range("N73").Value = -11849.223

$ws.Range("H92").Value = 36236.5
$ws.Range("J92").Value = 36236.5
$ws.Range("L92").Value = 36236.5
$ws.Range("N92").Value = -39980.5

$ws.Range("H97").Value = 428.2
$ws.Range("I97").Value = 399.1111
$ws.Range("J97").Value = 690
$ws.Range("K97").Value = 399.1111
$ws.Range("L97").Value = 690
$ws.Range("M97").Value = 96.88889999999998
$ws.Range("N97").Value = -1682

$ws.Range("H107").Value = 2371.9092
$ws.Range("I107").Value = 451
$ws.Range("J107").Value = 3972.6667
$ws.Range("K107").Value = 451
$ws.Range("L107").Value = 3972.6667
$ws.Range("M107").Value = 1469
$ws.Range("N107").Value = -7812.6667

$ws.Range("H132").Value = 1898.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2960
$ws.Range("I16").Value = 1933.3334
$ws.Range("K16").Value = 1933.3334
$ws.Range("M16").Value = -1763.3334

$ws.Range("H55").Value = 1713.2307
$ws.Range("I55").Value = 823.25
$ws.Range("J55").Value = 3137.2
$ws.Range("K55").Value = 823.25
$ws.Range("L55").Value = 3137.2
$ws.Range("M55").Value = -650.25
$ws.Range("N55").Value = -3483.2

$ws.Range("H61").Value = 2899.6667
$ws.Range("I61").Value = 2600
$ws.Range("J61").Value = 3049.5
$ws.Range("K61").Value = 2600
$ws.Range("L61").Value = 3049.5
$ws.Range("M61").Value = -2398
$ws.Range("N61").Value = -3453.5

$ws.Range("H100").Value = 3871.5715
$ws.Range("I100").Value = 3078.2222
$ws.Range("K100").Value = 3078.2222
$ws.Range("M100").Value = -2537.2222

$ws.Range("H113").Value = 2899.6667
$ws.Range("I113").Value = 2600
$ws.Range("J113").Value = 3049.5
$ws.Range("K113").Value = 2600
$ws.Range("L113").Value = 3049.5
$ws.Range("M113").Value = -430
$ws.Range("N113").Value = -7389.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = ""

$ws.Range("H45").Value = 20395.334
$ws.Range("J45").Value = 20395.334
$ws.Range("L45").Value = 20395.334
$ws.Range("N45").Value = -21377.334

$ws.Range("H56").Value = 50314
$ws.Range("J56").Value = 50314
$ws.Range("L56").Value = 50314
$ws.Range("N56").Value = -51742

$ws.Range("H100").Value = 8334697.5
$ws.Range("I100").Value = 8334697.5
$ws.Range("K100").Value = 16669395
$ws.Range("M100").Value = -16668854

$ws.Range("H107").Value = 1348.1666
$ws.Range("I107").Value = 1596.8572
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 4790.571599999999
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -2870.571599999999
$ws.Range("N107").Value = -6840

$ws.Range("H113").Value = 1245.3334
$ws.Range("I113").Value = 1379.6
$ws.Range("K113").Value = 4138.799999999999
$ws.Range("M113").Value = -1968.799999999999

$ws.Range("H132").Value = 7807.8184
$ws.Range("I132").Value = 7611.25
$ws.Range("J132").Value = 8332
$ws.Range("K132").Value = 22833.75
$ws.Range("L132").Value = 24996
$ws.Range("M132").Value = -20303.75
$ws.Range("N132").Value = -30056

$ws.Range("H136").Value = 3872.7307
$ws.Range("I136").Value = 3841.95
$ws.Range("K136").Value = 11525.85
$ws.Range("M136").Value = -8975.849999999999
